$d = $word.ActiveDocument

# 1. Update the subtitle text: "RequestSolved!" -> "E-commerce de joias e itens de artesanato"
$d.Content.Find.Execute("RequestSolved!", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "E-commerce de joias e itens de artesanato", 2)

# 2. Remove the first table (the "Tema" table) together with the three empty
#    paragraphs that separated it from the following table.
$temaTable = $d.Tables.Item(1)
$nextTable = $d.Tables.Item(2)

# Delete the blank paragraphs sitting between the two tables first (while the
# ranges are still easy to address), then drop the "Tema" table itself.
$gap = $d.Range($temaTable.Range.End, $nextTable.Range.Start)
$gap.Delete()

$d.Tables.Item(1).Delete()
